$wb = $excel.ActiveWorkbook

# --- "model" sheet: insert 8 new rows (session-variable / form-status
# bookkeeping fields) right after the existing "hh_head" row, and add a
# new "isSessionVariable" column used to flag rows assigned from the
# household form / session. ---
$model = $wb.Worksheets.Item("model")

# Insert 8 blank rows starting at row 5 (pushes everything else down).
$model.Range("A5:A12").EntireRow.Insert()

# New column header for the session-variable flag.
$model.Range("E1").Value = "isSessionVariable"

# New rows: name / type (/ isSessionVariable flag).
$model.Range("B5").Value = "id_candidate"
$model.Range("C5").Value = "string"
$model.Range("E5").Value = 1

$model.Range("B6").Value = "form_status_hh_member"
$model.Range("C6").Value = "integer"

$model.Range("B7").Value = "form_status_hh_member_absent"
$model.Range("C7").Value = "integer"

$model.Range("B8").Value = "form_status_hh_member_exit"
$model.Range("C8").Value = "integer"

$model.Range("B9").Value = "form_status_hh_member_new"
$model.Range("C9").Value = "integer"

$model.Range("B10").Value = "form_status_hh_member_questions"
$model.Range("C10").Value = "integer"

$model.Range("B11").Value = "form_status_hh_new_member_new"
$model.Range("C11").Value = "integer"

$model.Range("B12").Value = "form_status_hh_member_snake"
$model.Range("C12").Value = "integer"

# Make "model" the active sheet/tab, with the new column selected.
$model.Activate()
$model.Range("E1").Select()
